# Updated cryptos list with GitHub Actions scraped data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the Price column (D) cells we will touch as Text first, so that
# numeric-looking strings like "1.035" or "27.792.86" are kept verbatim
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.792.86"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.866.22"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("D4").Value = "1.035"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "323.98"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "1.032"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "0.4432"
$ws.Range("E7").Value = "  +1.70%  "
$ws.Range("D8").Value = "0.3815"
$ws.Range("E8").Value = "  +2.61%  "
$ws.Range("D9").Value = "0.07473"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").Value = "0.8904"
$ws.Range("E10").Value = "  +2.24%  "
$ws.Range("D11").Value = "21.78"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("D12").Value = "1.877.41"
$ws.Range("E12").Value = "  -4.85%  "
$ws.Range("D13").Value = "5.558"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "6.788"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "0.07209"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").Value = "84.54"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("D17").Value = "1.040"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "0.000009152"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").Value = "15.60"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "27.762.78"
$ws.Range("D22").Value = "5.326"
$ws.Range("E22").Value = "  +1.49%  "
$ws.Range("D23").Value = "11.34"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("D24").Value = "2.094.30"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").Value = "2.027"
$ws.Range("E25").Value = "  +6.32%  "
$ws.Range("D26").Value = "158.39"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").Value = "18.91"
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("D28").Value = "5.379"
$ws.Range("E28").Value = "  +2.65%  "
$ws.Range("E29").Value = "  +3.72%  "
$ws.Range("D30").Value = "119.25"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("D31").Value = "0.09046"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "1.241"
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").Value = "0.7824"
$ws.Range("E33").Value = "  +3.22%  "
$ws.Range("D34").Value = "3.023"
$ws.Range("E34").Value = "  +5.70%  "
$ws.Range("D35").Value = "4.613"
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").Value = "1.146"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").Value = "0.01992"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("D39").Value = "0.05365"
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("D40").Value = "2.893"
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("D41").Value = "0.5213"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("D42").Value = "0.1698"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("D43").Value = "6.925"
$ws.Range("E43").Value = "  +6.12%  "
$ws.Range("D44").Value = "8.748"
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("D45").Value = "111.43"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("D46").Value = "10.78"
$ws.Range("E46").Value = "  +2.92%  "
$ws.Range("D47").Value = "0.06677"
$ws.Range("E47").Value = "  +6.12%  "
$ws.Range("D48").Value = "1.036"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("D50").Value = "0.4739"
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("D51").Value = "1.921"
$ws.Range("E51").Value = "  +2.68%  "

# Restore the default (Normal) cell style on the Price column so the
# workbook styling matches the original (no custom number format left behind).
$ws.Range("D2:D51").Style = "Normal"

